$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Other"
$ws.Range("C1").Value = "+Other"
$ws.Range("D1").Value = "+I>C"

$ws.Range("D2").Select()
